# Updated cryptos list on Mon Jul  8 14:38:29 UTC 2024 with GitHub Actions
# Applies the latest scraped Price/Volume(1h) figures (and the OKB/FirstDigitalUSD
# and Filecoin/Mantle row-order swaps) to the crypto ranking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of cell -> new value updates taken from the latest scrape.
# Each entry is Column, Row, NewValue.
$updates = @(
    ,('D', 2, '55.588.49')
    ,('E', 2, '  -2.17%  ')
    ,('D', 3, '2.948.03')
    ,('E', 3, '  -0.67%  ')
    ,('D', 4, '1.01')
    ,('E', 4, '  +0.83%  ')
    ,('D', 5, '496.62')
    ,('E', 5, '  -0.18%  ')
    ,('D', 6, '135.83')
    ,('E', 6, '  -0.73%  ')
    ,('D', 7, '1.01')
    ,('E', 7, '  +0.61%  ')
    ,('D', 8, '0.427')
    ,('E', 8, '  +0.48%  ')
    ,('D', 9, '7.00')
    ,('E', 9, '  -4.15%  ')
    ,('D', 10, '0.105')
    ,('E', 10, '  -1.09%  ')
    ,('D', 11, '0.363')
    ,('E', 11, '  +2.07%  ')
    ,('D', 12, '3.516.94')
    ,('E', 12, '  +1.18%  ')
    ,('D', 13, '0.124')
    ,('E', 13, '  -2.62%  ')
    ,('D', 14, '25.89')
    ,('E', 14, '  +0.62%  ')
    ,('D', 15, '0.0000159')
    ,('E', 15, '  +1.75%  ')
    ,('D', 16, '56.128.83')
    ,('E', 16, '  -1.39%  ')
    ,('D', 17, '3.012.75')
    ,('E', 17, '  +1.87%  ')
    ,('D', 18, '5.89')
    ,('E', 18, '  -2.39%  ')
    ,('D', 19, '12.85')
    ,('E', 19, '  +2.15%  ')
    ,('D', 20, '7.77')
    ,('E', 20, '  -0.20%  ')
    ,('D', 21, '322.66')
    ,('E', 21, '  +1.47%  ')
    ,('D', 22, '0.996')
    ,('E', 22, '  -0.13%  ')
    ,('D', 23, '0.487')
    ,('E', 23, '  +0.54%  ')
    ,('D', 24, '64.29')
    ,('E', 24, '  +2.16%  ')
    ,('D', 25, '3.144.61')
    ,('E', 25, '  +1.94%  ')
    ,('D', 26, '1.01')
    ,('E', 26, '  +0.94%  ')
    ,('D', 27, '0.161')
    ,('E', 27, '  -0.24%  ')
    ,('D', 28, '0.0₃0865')
    ,('E', 28, '  -2.30%  ')
    ,('D', 29, '6.40')
    ,('E', 29, '  -1.81%  ')
    ,('D', 30, '6.90')
    ,('E', 30, '  -2.26%  ')
    ,('D', 31, '1.75')
    ,('E', 31, '  -0.08%  ')
    ,('D', 32, '1.15')
    ,('E', 32, '  -0.56%  ')
    ,('D', 33, '19.97')
    ,('E', 33, '  -0.58%  ')
    ,('D', 34, '151.55')
    ,('E', 34, '  -1.90%  ')
    ,('D', 35, '4.51')
    ,('E', 35, '  -1.71%  ')
    ,('D', 36, '5.71')
    ,('E', 36, '  +0.14%  ')
    ,('D', 37, '24.89')
    ,('E', 37, '  +4.08%  ')
    ,('D', 38, '1.23')
    ,('E', 38, '  -1.25%  ')
    ,('D', 39, '0.0652')
    ,('E', 39, '  -1.66%  ')
    ,('D', 40, '3.014.78')
    ,('E', 40, '  +0.49%  ')
    ,('B', 41, 'FirstDigitalUSD')
    ,('C', 41, 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd')
    ,('D', 41, '1.01')
    ,('E', 41, '  +1.06%  ')
    ,('B', 42, 'OKB')
    ,('C', 42, 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb')
    ,('D', 42, '36.48')
    ,('E', 42, '  -2.49%  ')
    ,('B', 43, 'Mantle')
    ,('C', 43, 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt')
    ,('D', 43, '0.650')
    ,('E', 43, '  +2.08%  ')
    ,('B', 44, 'Filecoin')
    ,('C', 44, 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil')
    ,('D', 44, '3.71')
    ,('E', 44, '  +0.26%  ')
    ,('D', 45, '2.147.58')
    ,('E', 45, '  -2.04%  ')
    ,('D', 46, '1.33')
    ,('E', 46, '  -3.02%  ')
    ,('D', 47, '0.926')
    ,('E', 47, '  -1.31%  ')
    ,('D', 48, '5.82')
    ,('E', 48, '  -1.54%  ')
    ,('D', 49, '0.0237')
    ,('E', 49, '  +1.07%  ')
    ,('D', 50, '19.51')
    ,('E', 50, '  +2.17%  ')
    ,('D', 51, '0.0844')
    ,('E', 51, '  -3.69%  ')
)

foreach ($u in $updates) {
    $col = $u[0]
    $row = $u[1]
    $newValue = $u[2]
    $cellRef = "$col$row"
    $cell = $ws.Range($cellRef)

    # The scraped values are always text (e.g. "7.00", "56.644.54", " +2.17% ").
    # Excel auto-converts plain numeric-looking strings (like "1.01" or "7.00")
    # to real numbers when assigned directly, which would lose the original
    # text formatting (trailing zeros, thousands "." separators, etc). Force
    # the cell to Text format first so the value is stored verbatim as a string,
    # then restore the default "Normal" style so no visible formatting changes.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}
